# Add KEGG_Pathways as a primary source.
#
# This touches three worksheets in the "Resource Interaction Table" workbook:
#   Sheet1 ("source interaction" table: source / predicate / target)
#     -> append a row: KEGG  -- has nomenclature mapping --> KEGG_Pathways
#   Sheet2 ("node / category" table)
#     -> append a row: KEGG_Pathways | Pathway
#   Sheet5 ("source interaction" table, same shape as Sheet1)
#     -> insert a row (keeping everything below shifted down by one):
#        KEGG -- has pathway link --> KEGG_Pathways

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: append new last row (row 95) for the KEGG / KEGG_Pathways
# nomenclature-mapping fact.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(95).Insert()
$ws1.Cells.Item(95, 1).Value = "KEGG"
$ws1.Cells.Item(95, 2).Value = "has nomenclature mapping"
$ws1.Cells.Item(95, 3).Value = "KEGG_Pathways"
$ws1.Rows.Item(95).RowHeight = 15.75

# ---------------------------------------------------------------------
# Sheet2: append new last row (row 45) registering the KEGG_Pathways
# node under the Pathway category.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows.Item(45).Insert()
$ws2.Cells.Item(45, 1).Value = "KEGG_Pathways"
$ws2.Cells.Item(45, 2).Value = "Pathway"
$ws2.Rows.Item(45).RowHeight = 15.75

# ---------------------------------------------------------------------
# Sheet5: insert a new row 41 (everything from the old row 41 onward
# shifts down by one) for the KEGG / KEGG_Pathways pathway-link fact.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Rows.Item(41).Insert()
$ws5.Cells.Item(41, 1).Value = "KEGG"
$ws5.Cells.Item(41, 2).Value = "has pathway link"
$ws5.Cells.Item(41, 3).Value = "KEGG_Pathways"
